# The observation rows 13-16 and 18-20 (all sharing the "Arvnäs, Ly lm"
# locality) were re-sorted/re-matched against their taxon records; row 17
# (Björktrast) was already correctly positioned and is left untouched.
# Net effect per row = a cyclic permutation of the whole-row content:
#   13<-14, 14<-15, 15<-18, 16<-19, 18<-20, 19<-13, 20<-16
# Below we only touch the individual cells whose value actually differs
# between "before" and "after" for each row.
#
# Note: plain numeric- or date-looking strings (e.g. "1", "2023-09-14")
# get auto-coerced to numbers/dates by COM's Range.Value setter, so for
# those specific cells we force the cell to Text format first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (<- old row 14): Sävsparv / Emberiza schoeniclus
$ws.Range("A13").Value = 112281199
$ws.Range("B13").Value = 57103
$ws.Range("E13").Value = 103057
$ws.Range("F13").Value = "Sävsparv"
$ws.Range("G13").Value = "Emberiza schoeniclus"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("I13").Value = ""

# Row 14 (<- old row 15): Rödvingetrast / Turdus iliacus
$ws.Range("A14").Value = 112279542
$ws.Range("B14").Value = 56841
$ws.Range("E14").Value = 103001
$ws.Range("F14").Value = "Rödvingetrast"
$ws.Range("G14").Value = "Turdus iliacus"
$ws.Range("H14").Value = "Linnaeus, 1766"

# Row 15 (<- old row 18): Lappsparv / Calcarius lapponicus
$ws.Range("A15").Value = 112281210
$ws.Range("B15").Value = 57076
$ws.Range("D15").Value = "VU"
$ws.Range("E15").Value = 103053
$ws.Range("F15").Value = "Lappsparv"
$ws.Range("G15").Value = "Calcarius lapponicus"
$ws.Range("H15").Value = "(Linnaeus, 1758)"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "1"
$ws.Range("M15").Value = "födosökande"

# Row 16 (<- old row 19): Spillkråka / Dryocopus martius
$ws.Range("A16").Value = 112292314
$ws.Range("B16").Value = 56446
$ws.Range("E16").Value = 100049
$ws.Range("F16").Value = "Spillkråka"
$ws.Range("G16").Value = "Dryocopus martius"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "1"
$ws.Range("M16").Value = "lockläte, övriga läten"
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2023-09-15"
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = "2023-09-15"

# Row 17 unchanged (Björktrast / Turdus pilaris)

# Row 18 (<- old row 20): Grönfink / Chloris chloris
$ws.Range("A18").Value = 112281233
$ws.Range("B18").Value = 57042
$ws.Range("D18").Value = "EN"
$ws.Range("E18").Value = 103042
$ws.Range("F18").Value = "Grönfink"
$ws.Range("G18").Value = "Chloris chloris"
$ws.Range("I18").Value = ""

# Row 19 (<- old row 13): Stenfalk / Falco columbarius
$ws.Range("A19").Value = 112279532
$ws.Range("B19").Value = 56476
$ws.Range("E19").Value = 102611
$ws.Range("F19").Value = "Stenfalk"
$ws.Range("G19").Value = "Falco columbarius"
$ws.Range("H19").Value = "Linnaeus, 1758"
$ws.Range("M19").Value = ""
$ws.Range("Y19").NumberFormat = "@"
$ws.Range("Y19").Value = "2023-09-14"
$ws.Range("AA19").NumberFormat = "@"
$ws.Range("AA19").Value = "2023-09-14"

# Row 20 (<- old row 16): Talltita / Poecile montanus
$ws.Range("A20").Value = 112281154
$ws.Range("B20").Value = 56575
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 103021
$ws.Range("F20").Value = "Talltita"
$ws.Range("G20").Value = "Poecile montanus"
$ws.Range("H20").Value = "(Conrad von Baldenstein, 1827)"
